$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet (tenant changed)
$ws.Name = "Brzoza"

# Update address line (house number and side changed)
$ws.Range("A2").Value = "Haferkamp  3 , 59192  Bergkamen,Doppelhaushälfte,re,EG,DG"

# Update tenant name
$ws.Range("B3").Value = "Brzoza"

# Row 5: Entwässerung -> Entwässerung Fläche, amount 608,82 -> 168,96
$ws.Range("A5").Value = "Entwässerung Fläche"
$ws.Range("C5").Value = "168,96"
$ws.Range("F5").Value = "168,96"

# Row 6: Entwässerung Fläche -> Entwässerung Verbrauch, amount 168,96 -> 784,02
$ws.Range("A6").Value = "Entwässerung Verbrauch"
$ws.Range("C6").Value = "784,02"
$ws.Range("F6").Value = "784,02"

# Row 8: Müllabfuhr amount 236,40 -> 118,20
$ws.Range("C8").Value = "118,20"
$ws.Range("F8").Value = "118,20"

# Row 12: Heizungswartung amount 122,68 -> 121,46
$ws.Range("C12").Value = "121,46"
$ws.Range("F12").Value = "121,46"

# Row 14: Summe 1.863,69 -> 1.919,47
$ws.Range("F14").Value = "1.919,47"
